$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of E3/F3 and E4/F4 so "Andre Barros-EAP" moves to column F
foreach ($r in 3, 4) {
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 5).Value2 = $fVal
    $ws.Cells.Item($r, 6).Value2 = $eVal
}
